$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Entities": mark SiteInfo (A13) and Tag (A17) as Done, add a new
# TagScore row (26), wire up its hyperlink, and touch the page setup /
# selection to mirror the authored edit.
# ---------------------------------------------------------------------------
$wsEntities = $wb.Worksheets.Item("Entities")

$wsEntities.Range("A13").Style = "Good"
$wsEntities.Range("A17").Style = "Good"

$wsEntities.Range("A26").Value = "TagScore"
$wsEntities.Range("A26").Style = "Good"

# Add the hyperlink first (Excel auto-creates a dedicated hyperlink style for
# the cell), then restyle it to match the plain "Hyperlink" look already used
# by the rest of column C by copying the format from an existing linked cell.
$wsEntities.Hyperlinks.Add($wsEntities.Range("C26"), "http://api.stackexchange.com/docs/types/tag-score")
$wsEntities.Range("C2").Copy()
$wsEntities.Range("C26").PasteSpecial(-4122)

$wsEntities.PageSetup.Orientation = 1

$wsEntities.Activate()
$wsEntities.Range("C26").Select()

# ---------------------------------------------------------------------------
# Sheet "Methods": mark the Search, Suggested Edits, Info, Tags and (part of)
# Users method rows as Done.
# ---------------------------------------------------------------------------
$wsMethods = $wb.Worksheets.Item("Methods")

$doneRows = @(40, 41, 43, 44, 46, 48, 49, 50, 51, 52, 53, 55, 56, 57, 58, 59, 60)
foreach ($r in $doneRows) {
    $wsMethods.Range("A$r").Value = "Done"
    $wsMethods.Range("A$r").Style = "Good"
}

$wsMethods.Activate()
$wsMethods.Range("A60").Select()
